# approx timing.xlsx — "Estimate with new moments" edit
#
# Adds a second timing-estimate block (columns D and F:L) to the
# "model_v2" sheet, alongside the existing emax-timing calculation in
# columns B:C. The new block mirrors the "From emax_timing code" numbers
# (Original / V2 / From sim_data / Parallel emax / V4 / V5 / V6 / v7
# scenarios) and a "Days" row under it. A few existing B:C values are
# also updated (grid size 1536->1800, emaX secs 58->69 captured in the
# new D column, N func evaluations 30->1000, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_v2")

# ---------------------------------------------------------------
# Row 3 / 4 headers
# ---------------------------------------------------------------
$ws.Range("B3").Value = "From sim_data"
$ws.Range("F3").Value = "From emax_timing code"

$ws.Range("C4").Value = "Original"
$ws.Range("D4").Value = "Parallel emax"

$ws.Range("G4").Value = "Original"
$ws.Range("H4").Value = "V2"
$ws.Range("I4").Value = "V4"
$ws.Range("J4").Value = "V5"
$ws.Range("K4").Value = "V6"
$ws.Range("L4").Value = "v7"
$ws.Range("G4:L4").Style = "Bold"

# ---------------------------------------------------------------
# Existing B:C block — updated values + new D column
# ---------------------------------------------------------------
$ws.Range("C5").Value = 1800
$ws.Range("D5").Value = 1800
$ws.Range("D5").Style = "Bold"

$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 50
$ws.Range("D6").Style = "Bold"

$ws.Range("D7").Value = 69

$ws.Range("C14").Value = 1000

# ---------------------------------------------------------------
# New "Cores for M parallel" / "Cores for emaxT parallel" /
# "Cores for emaxt parallel" rows (F5:L7)
# ---------------------------------------------------------------
$ws.Range("F5").Value = "Cores for M parallel"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 20
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 20
$ws.Range("L5").Value = 10
$ws.Range("G5:H5").Style = "Bold"

$ws.Range("F6").Value = "Cores for emaxT parallel"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 20
$ws.Range("L6").Value = 10
$ws.Range("G6:H6").Style = "Bold"

$ws.Range("F7").Value = "Cores for emaxt parallel"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 20
$ws.Range("L7").Value = 10
$ws.Range("G7:H7").Style = "Bold"

# ---------------------------------------------------------------
# Grid size / D / M rows (F8:L10)
# ---------------------------------------------------------------
$ws.Range("F8").Value = "Grid size"
$ws.Range("G8").Value = 1800
$ws.Range("H8").Value = 1800
$ws.Range("I8").Value = 1800
$ws.Range("J8").Value = 1800
$ws.Range("K8").Value = 1800
$ws.Range("L8").Value = 1800

$ws.Range("F9").Value = "D"
$ws.Range("G9").Value = 50
$ws.Range("H9").Value = 50
$ws.Range("I9").Value = 50
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 50
$ws.Range("L9").Value = 50

$ws.Range("F10").Value = "M"
$ws.Range("G10").Value = 1000
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1000

# ---------------------------------------------------------------
# Sample generation / Aux model generation / Total (secs) (F11:L13)
# ---------------------------------------------------------------
$ws.Range("F11").Value = "Sample generation (M datasets)"
$ws.Range("G11").Value = 228
$ws.Range("H11").Value = 106
$ws.Range("I11").Value = 85
$ws.Range("J11").Value = 85
$ws.Range("K11").Value = 73
$ws.Range("L11").Value = 58
$ws.Range("G11:L11").Style = "Bold"

$ws.Range("F12").Value = "Aux model generation"
$ws.Range("G12").Value = 0.000446
$ws.Range("H12").Value = 0.000446
$ws.Range("I12").Value = 0.000446
$ws.Range("J12").Value = 0.000446
$ws.Range("K12").Value = 0.000446
$ws.Range("L12").Value = 0.000446
$ws.Range("G12:L12").Style = "Bold"

$ws.Range("F13").Value = "Total (secs)"
$ws.Range("G13").Formula = "=SUM(G11:G12)"
$ws.Range("H13").Formula = "=SUM(H11:H12)"
$ws.Range("I13").Formula = "=SUM(I11:I12)"
$ws.Range("J13").Formula = "=SUM(J11:J12)"
$ws.Range("K13").Formula = "=SUM(K11:K12)"
$ws.Range("L13").Formula = "=SUM(L11:L12)"
$ws.Range("G13:L13").Style = "Bold"

# ---------------------------------------------------------------
# N func evaluations / Total (hours) / Days (F15:L17)
# ---------------------------------------------------------------
$ws.Range("F15").Value = "N func evaluations"
$ws.Range("G15").Value = 2456
$ws.Range("H15").Value = 2456
$ws.Range("I15").Value = 2456
$ws.Range("J15").Value = 2456
$ws.Range("K15").Value = 2456
$ws.Range("L15").Value = 2456

$ws.Range("F16").Value = "Total (hours)"
$ws.Range("F16").Style = "Bold"
$ws.Range("G16").Formula = "=G15*G13/(60*60)"
$ws.Range("H16").Formula = "=H15*H13/(60*60)"
$ws.Range("I16").Formula = "=I15*I13/(60*60)"
$ws.Range("J16").Formula = "=J15*J13/(60*60)"
$ws.Range("K16").Formula = "=K15*K13/(60*60)"
$ws.Range("L16").Formula = "=L15*L13/(60*60)"
$ws.Range("G16:L16").Style = "Bold"

$ws.Range("F17").Value = "Days"
$ws.Range("F17").Style = "Bold"
$ws.Range("G17").Formula = "=G16/24"
$ws.Range("H17").Formula = "=H16/24"
$ws.Range("I17").Formula = "=I16/24"
$ws.Range("J17").Formula = "=J16/24"
$ws.Range("K17").Formula = "=K16/24"
$ws.Range("L17").Formula = "=L16/24"
$ws.Range("G17:L17").HorizontalAlignment = -4108
$ws.Range("G17:L17").Font.Bold = $true
$ws.Range("G17:L17").NumberFormat = "0.0"

# ---------------------------------------------------------------
# Best-effort column widths for the new columns (engine quantises
# widths to 1/6-character steps, so these land as close as possible
# to the canonical bestFit pixel widths).
# ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 12.333333333333334
$ws.Columns.Item(6).ColumnWidth = 29.0
$ws.Columns.Item(7).ColumnWidth = 12.833333333333334
$ws.Columns.Item(8).ColumnWidth = 11.833333333333334

# ---------------------------------------------------------------
# Selection cursor ends on I17 (per the saved sheetView) and the page
# is explicitly set to portrait orientation (as in the target file).
# ---------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("I17").Select()
